$d = $word.ActiveDocument

$replacements = @(
    @("24÷7=", "67÷5="),
    @("70÷5=", "18÷3="),
    @("47÷5=", "52÷5="),
    @("98÷5=", "30÷3="),
    @("44÷8=", "87÷7="),
    @("96÷7=", "16÷2="),
    @("15÷3=", "86÷3="),
    @("26÷5=", "87÷3="),
    @("41÷7=", "57÷8="),
    @("48÷6=", "76÷4="),
    @("16÷8=", "17÷5="),
    @("80÷3=", "90÷8="),
    @("14÷7=", "59÷4="),
    @("17÷2=", "75÷3="),
    @("28÷2=", "13÷3="),
    @("53÷9=", "79÷9="),
    @("48÷7=", "96÷3="),
    @("54÷7=", "64÷4="),
    @("15÷4=", "72÷6="),
    @("58÷5=", "42÷9="),
    @("76÷8=", "39÷3="),
    @("35÷9=", "81÷9="),
    @("64÷5=", "24÷4="),
    @("73÷8=", "18÷9="),
    @("87÷2=", "22÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
